$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (old D -> E) to hold a "Tag" value
# derived from the project name in column C.
$ws.Columns("D").Insert()

$lastRow = 14

$ws.Range("D1").Value = "Tag"

for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 3).Value2
    $tag = $name -replace "[^A-Za-z0-9]+", "_"
    $tag = $tag.Trim("_")
    $ws.Cells.Item($r, 4).Value = $tag
}

# Copy column C's formatting onto the new column D (header + data rows)
$ws.Range("C1:C14").Copy()
$ws.Range("D1:D14").PasteSpecial(-4122)

# Match column widths to the newly laid-out table
$ws.Columns("D").ColumnWidth = 34
$ws.Columns("B").ColumnWidth = 20.42

[void]$ws.Range("D10").Select()
